# Append a new row (row 62) with sensor-reading data to each of the four
# worksheets, extending the used range from A1:I61 to A1:I62 on every sheet.

$wb = $excel.ActiveWorkbook

# Data for the new row on each worksheet, in sheet order.
$rows = @(
    @{
        Sheet = 1
        A = "2025-03-06 21:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "568631262647113770877196"
        H = 400
        I = 13
    },
    @{
        Sheet = 2
        A = "2025-03-06 21:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "568631262647113770942732"
        H = 400
        I = 14
    },
    @{
        Sheet = 3
        A = "2025-03-06 21:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 255
    },
    @{
        Sheet = 4
        A = "2025-03-06 21:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 3
    }
)

foreach ($rowInfo in $rows) {
    $ws = $wb.Worksheets.Item($rowInfo.Sheet)
    $newRow = 62

    $ws.Cells.Item($newRow, 1).Value2 = $rowInfo.A
    $ws.Cells.Item($newRow, 2).Value2 = $rowInfo.B
    $ws.Cells.Item($newRow, 3).Value2 = $rowInfo.C
    $ws.Cells.Item($newRow, 4).Value2 = $rowInfo.D
    $ws.Cells.Item($newRow, 5).Value2 = $rowInfo.E
    $ws.Cells.Item($newRow, 6).Value2 = $rowInfo.F

    # Column G holds a long digit string that must stay text (it exceeds
    # numeric precision), so force a text format before assigning it, then
    # restore the default "Normal" style so no stray formatting lingers.
    $gCell = $ws.Cells.Item($newRow, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value2 = $rowInfo.G
    $gCell.Style = "Normal"

    $ws.Cells.Item($newRow, 8).Value2 = $rowInfo.H
    $ws.Cells.Item($newRow, 9).Value2 = $rowInfo.I
}
